$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Through 2022-08-29"
$ws.Range("B1").Value = "August 2022 (through August 29)"

$ws.Range("B2").Value = 16
$ws.Range("J2").Value = 8
$ws.Range("AH3").Value = 2
$ws.Range("AX4").Value = 2
$ws.Range("AX5").Value = 10
$ws.Range("J7").Value = 7
$ws.Range("AX7").Value = 4
$ws.Range("BF7").Value = 2
$ws.Range("B10").Value = 2
$ws.Range("R26").Value = 2
$ws.Range("BF33").Value = 2
$ws.Range("AP34").Value = 2
$ws.Range("J35").Value = 1
$ws.Range("AP35").Value = 1
$ws.Range("AH36").Value = 1
$ws.Range("B50").Value = 4
$ws.Range("AH66").Value = 1
$ws.Range("AX74").Value = 2
